$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138, pushing existing rows 138:190 down to 139:191
$ws.Rows("138:138").Insert()

# Populate the newly inserted row 138 with the new data record
$ws.Range("A138").Value = 8
$ws.Range("B138").Value = "Terminal La Palmera de La Serena"
$ws.Range("C138").Value = "Coquimbo"
$ws.Range("D138").Value = 44917
$ws.Range("E138").Value = 4
$ws.Range("F138").Value = 100112040
$ws.Range("G138").Value = "Cilantro"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 2800
$ws.Range("K138").Value = 2500
$ws.Range("L138").Value = 3000
$ws.Range("M138").Value = 2750
$ws.Range("N138").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O138").Value = "Provincia del Elquí"
$ws.Range("P138").Value = 1833
$ws.Range("Q138").Value = 1.5
$ws.Range("R138").Value = "Hortaliza"

# Preserve the date cell style (numFmt) used by the rest of column D
$ws.Range("D138").NumberFormat = $ws.Range("D139").NumberFormat
